# January 12 / Make survey tool
# Adds a new "survey" (설문조사) table design section to Sheet1:
#   - 질문지 테이블 (questionnaire table, "servey") at rows 80-84
#   - 응답 테이블 (answer table) at rows 86-90

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Section header
$ws.Range("A78").Value = "설문조사"

# --- 질문지 테이블 (servey) ---
$ws.Range("A80").Value = "질문지 테이블"
$ws.Range("B80").Value = "servey"

$ws.Range("A81").Value = "타이틀"
$ws.Range("B81").Value = "일련번호"
$ws.Range("C81").Value = "질문"
$ws.Range("D81").Value = "선택내용1"
$ws.Range("E81").Value = "선택내용2"
$ws.Range("F81").Value = "선택내용3"
$ws.Range("G81").Value = "선택내용4"
$ws.Range("H81").Value = "서비스여부(0, 1)"

$ws.Range("A82").Value = "필드명"
$ws.Range("B82").Value = "no"
$ws.Range("C82").Value = "question"
$ws.Range("D82").Value = "select1"
$ws.Range("E82").Value = "select2"
$ws.Range("F82").Value = "select3"
$ws.Range("G82").Value = "select4"
$ws.Range("H82").Value = "status"

$ws.Range("A83").Value = "타입"
$ws.Range("B83").Value = "number"
$ws.Range("C83").Value = "varchar2(4000)"
$ws.Range("D83").Value = "varchar2(500)"
$ws.Range("E83").Value = "varchar2(500)"
$ws.Range("F83").Value = "varchar2(500)"
$ws.Range("G83").Value = "varchar2(500)"
$ws.Range("H83").Value = "char(1)"

$ws.Range("A84").Value = "제약조건"
$ws.Range("B84").Value = "not null"
$ws.Range("C84").Value = "not null"
$ws.Range("D84").Value = "not null"
$ws.Range("E84").Value = "not null"
$ws.Range("F84").Value = "not null"
$ws.Range("G84").Value = "not null"
$ws.Range("H84").Value = "default '0'"

# --- 응답 테이블 ---
$ws.Range("A86").Value = "응답 테이블"

$ws.Range("A87").Value = "타이틀"
$ws.Range("B87").Value = "일련번호"

$ws.Range("A88").Value = "필드명"
$ws.Range("B88").Value = "no"
$ws.Range("C88").Value = "servey_no"
$ws.Range("D88").Value = "servey_answer"

$ws.Range("A89").Value = "타입"
$ws.Range("B89").Value = "number"
$ws.Range("C89").Value = "number"
$ws.Range("D89").Value = "number"

$ws.Range("A90").Value = "제약조건"
$ws.Range("B90").Value = "not null"
$ws.Range("C90").Value = "not null"
$ws.Range("D90").Value = "not null"

# Update selection to match the final cursor position
$ws.Range("C91").Select()
